# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.187.08"
$ws.Range("E2").Value = "  +0.37%  "

$ws.Range("D3").Value = "1.602.10"

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("E5").Value = "  -0.11%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "303.17"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.38%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.3781"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "51.93"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +3.32%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.3624"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.92%  "

$ws.Range("E10").Value = "  -0.35%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.09%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.08115"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.32%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "22.82"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.17%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.601"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.21%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.417"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.24%  "

$ws.Range("E16").Value = "  -1.25%  "

$ws.Range("D17").Value = "1.603.62"
$ws.Range("E17").Value = "  +0.09%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "93.88"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.94%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06871"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.14%  "

$ws.Range("E20").Value = "  -1.21%  "

$ws.Range("E21").Value = "  -0.81%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.08%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "12.97"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -1.20%  "

$ws.Range("D24").Value = "23.192.12"
$ws.Range("E24").Value = "  +0.39%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "3.023"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +7.87%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.390"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +1.40%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "21.22"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.22%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "150.02"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.32%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "5.247"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.41%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "133.86"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.19%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "2.365"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.67%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "6.762"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.22%  "

$ws.Range("D33").Value = "1.781.26"
$ws.Range("E33").Value = "  +0.16%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.9675"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.89%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.07516"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -2.40%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.02727"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.19%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "10.23"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.76%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.2523"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.16%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.08797"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -1.22%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "6.090"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -3.01%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.7111"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.26%  "

$ws.Range("E42").Value = "  -0.42%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "12.63"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.05%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "15.64"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +1.23%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.6553"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.32%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.312"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.53%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "4.018"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.34%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "132.36"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.19%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.07947"
$c.Style = "Normal"

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.205"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -2.92%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.208"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.26%  "
